$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Password" / "Role" header columns (bold, like the other header cells) ---
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Role"
$ws.Range("D1:E1").Font.Bold = $true

# --- Staff rows: password + role values ---
$ws.Range("D2").Value = "password"
$ws.Range("E2").Value = "Staff"

$ws.Range("D3").Value = "password"
$ws.Range("E3").Value = "Staff"

$ws.Range("D4").Value = "password"
$ws.Range("E4").Value = "Staff"

$ws.Range("D5").Value = "password"
$ws.Range("E5").Value = "Staff"

# Row 6 ("Arvind") keeps "1" as text (not a number) in the Password column
$ws.Range("D6").Value = "'1"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "Staff"

# --- Fix trailing space in the staff name "Arvind " -> "Arvind" ---
$ws.Range("A6").Value = "Arvind"

$ws.Range("E6").Select()
